$d = $word.ActiveDocument

# Locate the "D - 0-28 points" paragraph (grading scale line) which is
# immediately followed by a run of empty paragraphs before the closing
# "Project maintained at ..." line.
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("0-28 points")) {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'D - 0-28 points' paragraph"
}

# Remove 4 of the extra blank paragraphs that follow it (leftover filler
# the author forgot to trim), leaving the rest of the spacing intact.
$toRemove = 4
$firstIndex = $anchorIndex + 1
$lastIndex = $anchorIndex + $toRemove

$firstPara = $d.Paragraphs.Item($firstIndex)
$lastPara = $d.Paragraphs.Item($lastIndex)

$rng = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$rng.Delete()
